$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) verifyTheSearch: trim the two trailing blank rows (A1:B4 -> A1:A2)
# ---------------------------------------------------------------------------
$wsSearch = $wb.Worksheets.Item("verifyTheSearch")
$wsSearch.Range("A3:B4").Clear()

# ---------------------------------------------------------------------------
# 2) Append five new test-data sheets at the end of the workbook
# ---------------------------------------------------------------------------

# verifyDynamicDropdown
$prev = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $prev)
$ws.Name = "verifyDynamicDropdown"
$ws.Range("A1").Value = "AskOsmoseSearch"
$ws.Range("A2").Value = "Ind"
[void]$ws.Range("A2").Select()

# verifySearchResultsOnTyping
$prev = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $prev)
$ws.Name = "verifySearchResultsOnTyping"
$ws.Range("A1").Value = "TypeJointUse"
$ws.Range("A2").Value = "Joint Use"
[void]$ws.Range("A1:A2").Select()

# verifyClearButtonOnTyping
$prev = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $prev)
$ws.Name = "verifyClearButtonOnTyping"
$ws.Range("A1").Value = "TypeJointUseInSearchbox"
$ws.Range("A2").Value = "Joint Use"
[void]$ws.Range("H12").Select()

# verifyNoResultsFoundMessage
$prev = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $prev)
$ws.Name = "verifyNoResultsFoundMessage"
$ws.Range("A1").Value = "NoResults"
$ws.Range("A2").Value = "Noresultsfound"
[void]$ws.Range("D8").Select()

# verifySearchTextUpdation
$prev = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $prev)
$ws.Name = "verifySearchTextUpdation"
$ws.Range("A1").Value = "TypeJointUseTextUpdation"
$ws.Range("B1").Value = "TypePoleReplacementTextUpdation"
$ws.Range("A2").Value = "Joint Use"
$ws.Range("B2").Value = "Pole Replacement"
[void]$ws.Range("B1").Select()
